$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove stray border-only placeholder cells that no longer belong ---
$ws.Range("C4:K4").Clear()
$ws.Range("I6:K6").Clear()
$ws.Range("C8:K8").Clear()
$ws.Range("I10:K10").Clear()

# Row 12: C12:I12 keep their text but lose their special border formatting;
# J12:K12 were empty placeholder cells and are removed entirely.
$ws.Range("C12:I12").ClearFormats()
$ws.Range("J12:K12").Clear()

# --- Insert a new row above the old row 15 ("pair_style smatb ...") ---
# This pushes the old row15/16 down to 16/17 (and 18-24 down to 19-25),
# and Excel automatically re-points the C15/C16 references used by the
# F13/G13 formulas to C16/C17.
$ws.Rows("15:15").Insert()

# --- Populate the newly freed row 14 with the "check" helper row ---
$ws.Range("B14").Value = "check"
$ws.Range("C14").Formula = "=E13/((E5/C5+E9/C9)/2)"
$ws.Range("C14").NumberFormat = "0.0000"

$ws.Range("D14").Value = '<- using "R0 = q(A-B)/((q(A)/R0(A)+q(B)/R0(B))/2)"'
$boldPart = $ws.Range("D14").Characters(16, 35)
$boldPart.Font.Bold = $true

# Draw a thin top+bottom border across B14:H14, with closing left/right
# edges only at the two ends of the range (matching the original box).
# Border indices: 7=left, 8=top, 9=bottom, 10=right.
$ws.Range("B14:H14").Borders.Item(8).LineStyle = 1
$ws.Range("B14:H14").Borders.Item(9).LineStyle = 1
$ws.Range("B14").Borders.Item(7).LineStyle = 1
$ws.Range("H14").Borders.Item(10).LineStyle = 1

# --- Column H needs to be a bit wider to show the new note text ---
$ws.Columns("H:H").ColumnWidth = 8.8

# --- Restore the selected cell as recorded in the saved view state ---
$ws.Range("E16").Select()
